$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164, pushing existing row 164 (and below) down to 165+
$ws.Rows.Item(164).EntireRow.Insert()

# Populate the newly inserted row 164 with the new weekly record
$ws.Cells.Item(164, 1).Value = 5
$ws.Cells.Item(164, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(164, 3).Value = "Maule"
$ws.Cells.Item(164, 4).Value = 45029
$ws.Cells.Item(164, 4).NumberFormat = $ws.Cells.Item(165, 4).NumberFormat
$ws.Cells.Item(164, 5).Value = 7
$ws.Cells.Item(164, 6).Value = 100112031
$ws.Cells.Item(164, 7).Value = "Poroto verde"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 150
$ws.Cells.Item(164, 11).Value = 25000
$ws.Cells.Item(164, 12).Value = 25000
$ws.Cells.Item(164, 13).Value = 25000
$ws.Cells.Item(164, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(164, 15).Value = "Región del Maule"
$ws.Cells.Item(164, 16).Value = 1000
$ws.Cells.Item(164, 17).Value = 25
$ws.Cells.Item(164, 18).Value = "Hortaliza"
